# Actualización automática 2025-10-15 14:30:09
#
# Updates the October sales figures for cliente "MANCHENO PINO HERVIN
# SANTIAGO" (asesor "ALMEIDA CUATIN JHONATHANN CARLOS"):
#   - INODOROS    sales: 71.09999999999999 -> 173.7   (+102.6)
#   - PORCELANATO sales: 2680.32           -> 2974.62 (+294.3)
#
# and propagates the resulting totals/derived figures across the other
# two sheets of the workbook ("VENTA MENSUAL" and "CUMPLIMIENTO MENSUAL").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" - raw per-client / per-group sales figures
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("H20").Value = 173.7      # INODOROS
$wsGrupo.Range("M20").Value = 2974.62    # PORCELANATO

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" - monthly totals per client, plus grand totals
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F20").Value = 3532.65   # octubre, same client as above
$wsMensual.Range("F36").Value = 14209.13  # octubre grand total

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL" - budget vs. actual-sale compliance
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 6: INODOROS
$wsCumpl.Range("D6").Value = 173.7
$wsCumpl.Range("E6").Value = 640.423430808873
$wsCumpl.Range("F6").Value = 0.2133583108244658

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 12795.72
$wsCumpl.Range("E12").Value = 8905.550000000001
$wsCumpl.Range("F12").Value = 0.5896300078290349

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 14209.13
$wsCumpl.Range("E14").Value = 22376.43723718182
$wsCumpl.Range("F14").Value = 0.3883807488314489

# Column E widened slightly (22 -> 23) as a side effect of the refresh.
# (The stored OOXML <col width> is ColumnWidth + 0.8333333333333334, the
# standard Excel character-to-width padding offset, so back it out here
# to land exactly on width="23".)
$wsCumpl.Columns.Item(5).ColumnWidth = 22.166666666666668
